$wb = $excel.ActiveWorkbook

# Insert a new "State" column into the hotel_info sheet, between
# "Hotel_Name" and "City", with value "Louisiana" for the existing row.
$hotelWs = $wb.Worksheets.Item("hotel_info")
$hotelWs.Columns.Item(3).Insert()
$hotelWs.Cells.Item(1, 3).Value = "State"
$hotelWs.Cells.Item(2, 3).Value = "Louisiana"

# Reorder the sheet tabs so "review_info" becomes the first sheet and
# "hotel_info" becomes the second sheet.
$reviewWs = $wb.Worksheets.Item("review_info")
$reviewWs.Move($wb.Worksheets.Item(1))
